$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 45-132 (shifted record data) ---
$ws.Cells.Item(45, 4).Value2 = 44811
$ws.Cells.Item(45, 10).Value2 = 150
$ws.Cells.Item(45, 11).Value2 = 22000
$ws.Cells.Item(45, 12).Value2 = 23000
$ws.Cells.Item(45, 13).Value2 = 22500
$ws.Cells.Item(45, 16).Value2 = 1125
$ws.Cells.Item(46, 4).Value2 = 44811
$ws.Cells.Item(46, 9).Value2 = "Segunda"
$ws.Cells.Item(46, 10).Value2 = 160
$ws.Cells.Item(46, 11).Value2 = 18000
$ws.Cells.Item(46, 12).Value2 = 19000
$ws.Cells.Item(46, 13).Value2 = 18500
$ws.Cells.Item(46, 16).Value2 = 925
$ws.Cells.Item(47, 4).Value2 = 44305
$ws.Cells.Item(47, 10).Value2 = 160
$ws.Cells.Item(47, 11).Value2 = 20000
$ws.Cells.Item(47, 12).Value2 = 21000
$ws.Cells.Item(47, 13).Value2 = 20500
$ws.Cells.Item(47, 16).Value2 = 1025
$ws.Cells.Item(48, 4).Value2 = 44750
$ws.Cells.Item(48, 10).Value2 = 140
$ws.Cells.Item(48, 11).Value2 = 33000
$ws.Cells.Item(48, 12).Value2 = 35000
$ws.Cells.Item(48, 13).Value2 = 34000
$ws.Cells.Item(48, 16).Value2 = 1700
$ws.Cells.Item(49, 4).Value2 = 44610
$ws.Cells.Item(49, 10).Value2 = 120
$ws.Cells.Item(49, 11).Value2 = 37000
$ws.Cells.Item(49, 12).Value2 = 38000
$ws.Cells.Item(49, 13).Value2 = 37500
$ws.Cells.Item(49, 16).Value2 = 1875
$ws.Cells.Item(50, 4).Value2 = 44468
$ws.Cells.Item(50, 9).Value2 = "Primera"
$ws.Cells.Item(50, 10).Value2 = 120
$ws.Cells.Item(50, 11).Value2 = 27000
$ws.Cells.Item(50, 12).Value2 = 28000
$ws.Cells.Item(50, 13).Value2 = 27500
$ws.Cells.Item(50, 16).Value2 = 1375
$ws.Cells.Item(51, 4).Value2 = 44372
$ws.Cells.Item(51, 9).Value2 = "Primera"
$ws.Cells.Item(51, 10).Value2 = 370
$ws.Cells.Item(51, 12).Value2 = 19000
$ws.Cells.Item(51, 13).Value2 = 18459
$ws.Cells.Item(51, 16).Value2 = 923
$ws.Cells.Item(52, 4).Value2 = 44372
$ws.Cells.Item(52, 9).Value2 = "Segunda"
$ws.Cells.Item(52, 10).Value2 = 150
$ws.Cells.Item(52, 11).Value2 = 14000
$ws.Cells.Item(52, 12).Value2 = 15000
$ws.Cells.Item(52, 13).Value2 = 14333
$ws.Cells.Item(52, 16).Value2 = 717
$ws.Cells.Item(53, 4).Value2 = 44722
$ws.Cells.Item(53, 9).Value2 = "Segunda"
$ws.Cells.Item(53, 10).Value2 = 150
$ws.Cells.Item(53, 11).Value2 = 18000
$ws.Cells.Item(53, 12).Value2 = 20000
$ws.Cells.Item(53, 13).Value2 = 19000
$ws.Cells.Item(53, 16).Value2 = 950
$ws.Cells.Item(54, 4).Value2 = 44312
$ws.Cells.Item(54, 9).Value2 = "Primera"
$ws.Cells.Item(54, 11).Value2 = 20000
$ws.Cells.Item(54, 12).Value2 = 21000
$ws.Cells.Item(54, 13).Value2 = 20500
$ws.Cells.Item(54, 16).Value2 = 1025
$ws.Cells.Item(55, 4).Value2 = 44281
$ws.Cells.Item(55, 10).Value2 = 100
$ws.Cells.Item(55, 11).Value2 = 29000
$ws.Cells.Item(55, 12).Value2 = 30000
$ws.Cells.Item(55, 13).Value2 = 29500
$ws.Cells.Item(55, 16).Value2 = 1475
$ws.Cells.Item(56, 4).Value2 = 44281
$ws.Cells.Item(56, 10).Value2 = 120
$ws.Cells.Item(56, 11).Value2 = 24000
$ws.Cells.Item(56, 12).Value2 = 25000
$ws.Cells.Item(56, 13).Value2 = 24500
$ws.Cells.Item(56, 16).Value2 = 1225
$ws.Cells.Item(57, 4).Value2 = 44428
$ws.Cells.Item(58, 4).Value2 = 44428
$ws.Cells.Item(58, 9).Value2 = "Segunda"
$ws.Cells.Item(58, 10).Value2 = 160
$ws.Cells.Item(58, 11).Value2 = 17000
$ws.Cells.Item(58, 12).Value2 = 18000
$ws.Cells.Item(58, 13).Value2 = 17500
$ws.Cells.Item(58, 16).Value2 = 875
$ws.Cells.Item(59, 4).Value2 = 44799
$ws.Cells.Item(59, 9).Value2 = "Primera"
$ws.Cells.Item(59, 10).Value2 = 160
$ws.Cells.Item(59, 11).Value2 = 19000
$ws.Cells.Item(59, 12).Value2 = 20000
$ws.Cells.Item(59, 13).Value2 = 19500
$ws.Cells.Item(59, 16).Value2 = 975
$ws.Cells.Item(60, 4).Value2 = 44792
$ws.Cells.Item(60, 10).Value2 = 140
$ws.Cells.Item(60, 11).Value2 = 20000
$ws.Cells.Item(60, 12).Value2 = 21000
$ws.Cells.Item(60, 13).Value2 = 20500
$ws.Cells.Item(60, 16).Value2 = 1025
$ws.Cells.Item(61, 4).Value2 = 44792
$ws.Cells.Item(61, 10).Value2 = 150
$ws.Cells.Item(61, 11).Value2 = 15000
$ws.Cells.Item(61, 12).Value2 = 16000
$ws.Cells.Item(61, 13).Value2 = 15500
$ws.Cells.Item(61, 16).Value2 = 775
$ws.Cells.Item(62, 4).Value2 = 44424
$ws.Cells.Item(62, 10).Value2 = 130
$ws.Cells.Item(62, 11).Value2 = 23000
$ws.Cells.Item(62, 12).Value2 = 24000
$ws.Cells.Item(62, 13).Value2 = 23500
$ws.Cells.Item(62, 16).Value2 = 1175
$ws.Cells.Item(63, 4).Value2 = 44424
$ws.Cells.Item(63, 10).Value2 = 140
$ws.Cells.Item(63, 11).Value2 = 20000
$ws.Cells.Item(63, 12).Value2 = 21000
$ws.Cells.Item(63, 13).Value2 = 20500
$ws.Cells.Item(63, 16).Value2 = 1025
$ws.Cells.Item(64, 4).Value2 = 44568
$ws.Cells.Item(64, 11).Value2 = 14000
$ws.Cells.Item(64, 12).Value2 = 15000
$ws.Cells.Item(64, 13).Value2 = 14500
$ws.Cells.Item(64, 16).Value2 = 725
$ws.Cells.Item(65, 4).Value2 = 44568
$ws.Cells.Item(65, 9).Value2 = "Segunda"
$ws.Cells.Item(65, 11).Value2 = 10000
$ws.Cells.Item(65, 12).Value2 = 11000
$ws.Cells.Item(65, 13).Value2 = 10500
$ws.Cells.Item(65, 16).Value2 = 525
$ws.Cells.Item(66, 4).Value2 = 44308
$ws.Cells.Item(66, 9).Value2 = "Primera"
$ws.Cells.Item(66, 11).Value2 = 19000
$ws.Cells.Item(66, 12).Value2 = 20000
$ws.Cells.Item(66, 13).Value2 = 19500
$ws.Cells.Item(66, 16).Value2 = 975
$ws.Cells.Item(67, 4).Value2 = 44323
$ws.Cells.Item(67, 10).Value2 = 120
$ws.Cells.Item(67, 11).Value2 = 21000
$ws.Cells.Item(67, 12).Value2 = 22000
$ws.Cells.Item(67, 13).Value2 = 21500
$ws.Cells.Item(67, 16).Value2 = 1075
$ws.Cells.Item(68, 4).Value2 = 44323
$ws.Cells.Item(68, 10).Value2 = 120
$ws.Cells.Item(68, 11).Value2 = 18000
$ws.Cells.Item(68, 12).Value2 = 19000
$ws.Cells.Item(68, 13).Value2 = 18500
$ws.Cells.Item(68, 16).Value2 = 925
$ws.Cells.Item(69, 4).Value2 = 44477
$ws.Cells.Item(69, 10).Value2 = 130
$ws.Cells.Item(69, 11).Value2 = 26000
$ws.Cells.Item(69, 12).Value2 = 27000
$ws.Cells.Item(69, 13).Value2 = 26500
$ws.Cells.Item(69, 16).Value2 = 1325
$ws.Cells.Item(70, 4).Value2 = 44477
$ws.Cells.Item(70, 10).Value2 = 140
$ws.Cells.Item(70, 11).Value2 = 23000
$ws.Cells.Item(70, 12).Value2 = 24000
$ws.Cells.Item(70, 13).Value2 = 23500
$ws.Cells.Item(70, 16).Value2 = 1175
$ws.Cells.Item(71, 4).Value2 = 44536
$ws.Cells.Item(71, 11).Value2 = 15000
$ws.Cells.Item(71, 12).Value2 = 16000
$ws.Cells.Item(71, 13).Value2 = 15500
$ws.Cells.Item(71, 16).Value2 = 775
$ws.Cells.Item(72, 4).Value2 = 44536
$ws.Cells.Item(72, 9).Value2 = "Segunda"
$ws.Cells.Item(72, 10).Value2 = 120
$ws.Cells.Item(72, 11).Value2 = 13000
$ws.Cells.Item(72, 12).Value2 = 14000
$ws.Cells.Item(72, 13).Value2 = 13500
$ws.Cells.Item(72, 16).Value2 = 675
$ws.Cells.Item(73, 4).Value2 = 44403
$ws.Cells.Item(73, 9).Value2 = "Primera"
$ws.Cells.Item(73, 11).Value2 = 29000
$ws.Cells.Item(73, 12).Value2 = 30000
$ws.Cells.Item(73, 13).Value2 = 29500
$ws.Cells.Item(73, 16).Value2 = 1475
$ws.Cells.Item(74, 4).Value2 = 44407
$ws.Cells.Item(74, 10).Value2 = 140
$ws.Cells.Item(74, 11).Value2 = 33000
$ws.Cells.Item(74, 12).Value2 = 34000
$ws.Cells.Item(74, 13).Value2 = 33500
$ws.Cells.Item(74, 16).Value2 = 1675
$ws.Cells.Item(75, 4).Value2 = 44407
$ws.Cells.Item(75, 9).Value2 = "Segunda"
$ws.Cells.Item(75, 10).Value2 = 120
$ws.Cells.Item(75, 11).Value2 = 30000
$ws.Cells.Item(75, 12).Value2 = 31000
$ws.Cells.Item(75, 13).Value2 = 30500
$ws.Cells.Item(75, 16).Value2 = 1525
$ws.Cells.Item(76, 4).Value2 = 44169
$ws.Cells.Item(76, 9).Value2 = "Primera"
$ws.Cells.Item(76, 10).Value2 = 160
$ws.Cells.Item(76, 11).Value2 = 18000
$ws.Cells.Item(76, 12).Value2 = 20000
$ws.Cells.Item(76, 13).Value2 = 19000
$ws.Cells.Item(76, 16).Value2 = 950
$ws.Cells.Item(77, 4).Value2 = 44790
$ws.Cells.Item(77, 10).Value2 = 140
$ws.Cells.Item(78, 4).Value2 = 44790
$ws.Cells.Item(78, 10).Value2 = 150
$ws.Cells.Item(78, 11).Value2 = 21000
$ws.Cells.Item(78, 12).Value2 = 22000
$ws.Cells.Item(78, 13).Value2 = 21500
$ws.Cells.Item(78, 16).Value2 = 1075
$ws.Cells.Item(79, 4).Value2 = 44785
$ws.Cells.Item(79, 10).Value2 = 130
$ws.Cells.Item(79, 11).Value2 = 24000
$ws.Cells.Item(79, 12).Value2 = 25000
$ws.Cells.Item(79, 13).Value2 = 24500
$ws.Cells.Item(79, 16).Value2 = 1225
$ws.Cells.Item(80, 4).Value2 = 44785
$ws.Cells.Item(80, 9).Value2 = "Segunda"
$ws.Cells.Item(80, 10).Value2 = 120
$ws.Cells.Item(80, 11).Value2 = 19000
$ws.Cells.Item(80, 12).Value2 = 20000
$ws.Cells.Item(80, 13).Value2 = 19500
$ws.Cells.Item(80, 16).Value2 = 975
$ws.Cells.Item(81, 4).Value2 = 44554
$ws.Cells.Item(81, 9).Value2 = "Primera"
$ws.Cells.Item(81, 10).Value2 = 160
$ws.Cells.Item(81, 11).Value2 = 10000
$ws.Cells.Item(81, 12).Value2 = 11000
$ws.Cells.Item(81, 13).Value2 = 10500
$ws.Cells.Item(81, 16).Value2 = 525
$ws.Cells.Item(82, 4).Value2 = 44757
$ws.Cells.Item(82, 10).Value2 = 150
$ws.Cells.Item(82, 11).Value2 = 33000
$ws.Cells.Item(82, 12).Value2 = 35000
$ws.Cells.Item(82, 13).Value2 = 34000
$ws.Cells.Item(82, 16).Value2 = 1700
$ws.Cells.Item(83, 4).Value2 = 44757
$ws.Cells.Item(83, 9).Value2 = "Segunda"
$ws.Cells.Item(83, 10).Value2 = 170
$ws.Cells.Item(83, 11).Value2 = 30000
$ws.Cells.Item(83, 12).Value2 = 32000
$ws.Cells.Item(83, 13).Value2 = 31000
$ws.Cells.Item(83, 16).Value2 = 1550
$ws.Cells.Item(84, 4).Value2 = 44540
$ws.Cells.Item(84, 10).Value2 = 160
$ws.Cells.Item(84, 11).Value2 = 12000
$ws.Cells.Item(84, 12).Value2 = 13000
$ws.Cells.Item(84, 13).Value2 = 12500
$ws.Cells.Item(84, 16).Value2 = 625
$ws.Cells.Item(85, 4).Value2 = 44298
$ws.Cells.Item(85, 10).Value2 = 140
$ws.Cells.Item(85, 11).Value2 = 20000
$ws.Cells.Item(85, 12).Value2 = 21000
$ws.Cells.Item(85, 13).Value2 = 20500
$ws.Cells.Item(85, 16).Value2 = 1025
$ws.Cells.Item(86, 4).Value2 = 44498
$ws.Cells.Item(86, 9).Value2 = "Primera"
$ws.Cells.Item(86, 10).Value2 = 130
$ws.Cells.Item(86, 11).Value2 = 37000
$ws.Cells.Item(86, 12).Value2 = 38000
$ws.Cells.Item(86, 13).Value2 = 37500
$ws.Cells.Item(86, 16).Value2 = 1875
$ws.Cells.Item(87, 4).Value2 = 44410
$ws.Cells.Item(87, 11).Value2 = 29000
$ws.Cells.Item(87, 12).Value2 = 30000
$ws.Cells.Item(87, 13).Value2 = 29500
$ws.Cells.Item(87, 16).Value2 = 1475
$ws.Cells.Item(88, 4).Value2 = 44410
$ws.Cells.Item(88, 10).Value2 = 120
$ws.Cells.Item(88, 11).Value2 = 27000
$ws.Cells.Item(88, 12).Value2 = 28000
$ws.Cells.Item(88, 13).Value2 = 27500
$ws.Cells.Item(88, 16).Value2 = 1375
$ws.Cells.Item(89, 4).Value2 = 44806
$ws.Cells.Item(89, 10).Value2 = 120
$ws.Cells.Item(89, 11).Value2 = 19000
$ws.Cells.Item(89, 12).Value2 = 20000
$ws.Cells.Item(89, 13).Value2 = 19500
$ws.Cells.Item(89, 16).Value2 = 975
$ws.Cells.Item(90, 4).Value2 = 44806
$ws.Cells.Item(90, 10).Value2 = 130
$ws.Cells.Item(90, 11).Value2 = 17000
$ws.Cells.Item(90, 12).Value2 = 18000
$ws.Cells.Item(90, 13).Value2 = 17500
$ws.Cells.Item(90, 16).Value2 = 875
$ws.Cells.Item(91, 4).Value2 = 44414
$ws.Cells.Item(91, 10).Value2 = 160
$ws.Cells.Item(91, 11).Value2 = 24000
$ws.Cells.Item(91, 12).Value2 = 25000
$ws.Cells.Item(91, 13).Value2 = 24500
$ws.Cells.Item(91, 16).Value2 = 1225
$ws.Cells.Item(92, 4).Value2 = 44414
$ws.Cells.Item(92, 10).Value2 = 140
$ws.Cells.Item(92, 11).Value2 = 21000
$ws.Cells.Item(92, 12).Value2 = 22000
$ws.Cells.Item(92, 13).Value2 = 21500
$ws.Cells.Item(92, 16).Value2 = 1075
$ws.Cells.Item(93, 4).Value2 = 44715
$ws.Cells.Item(93, 11).Value2 = 35000
$ws.Cells.Item(93, 12).Value2 = 38000
$ws.Cells.Item(93, 13).Value2 = 36500
$ws.Cells.Item(93, 16).Value2 = 1825
$ws.Cells.Item(94, 4).Value2 = 44715
$ws.Cells.Item(94, 11).Value2 = 30000
$ws.Cells.Item(94, 12).Value2 = 33000
$ws.Cells.Item(94, 13).Value2 = 31500
$ws.Cells.Item(94, 16).Value2 = 1575
$ws.Cells.Item(95, 4).Value2 = 44302
$ws.Cells.Item(95, 10).Value2 = 120
$ws.Cells.Item(95, 11).Value2 = 28000
$ws.Cells.Item(95, 12).Value2 = 30000
$ws.Cells.Item(95, 13).Value2 = 29000
$ws.Cells.Item(95, 16).Value2 = 1450
$ws.Cells.Item(96, 4).Value2 = 44302
$ws.Cells.Item(96, 9).Value2 = "Segunda"
$ws.Cells.Item(96, 10).Value2 = 120
$ws.Cells.Item(96, 11).Value2 = 19000
$ws.Cells.Item(96, 12).Value2 = 20000
$ws.Cells.Item(96, 13).Value2 = 19500
$ws.Cells.Item(96, 16).Value2 = 975
$ws.Cells.Item(97, 4).Value2 = 44418
$ws.Cells.Item(97, 10).Value2 = 130
$ws.Cells.Item(98, 4).Value2 = 44809
$ws.Cells.Item(98, 9).Value2 = "Primera"
$ws.Cells.Item(98, 10).Value2 = 140
$ws.Cells.Item(99, 4).Value2 = 44316
$ws.Cells.Item(100, 4).Value2 = 44316
$ws.Cells.Item(100, 9).Value2 = "Segunda"
$ws.Cells.Item(100, 11).Value2 = 20000
$ws.Cells.Item(100, 12).Value2 = 21000
$ws.Cells.Item(100, 13).Value2 = 20500
$ws.Cells.Item(100, 16).Value2 = 1025
$ws.Cells.Item(101, 4).Value2 = 44764
$ws.Cells.Item(101, 9).Value2 = "Primera"
$ws.Cells.Item(101, 11).Value2 = 24000
$ws.Cells.Item(101, 12).Value2 = 25000
$ws.Cells.Item(101, 13).Value2 = 24500
$ws.Cells.Item(101, 16).Value2 = 1225
$ws.Cells.Item(102, 4).Value2 = 44557
$ws.Cells.Item(102, 11).Value2 = 13000
$ws.Cells.Item(102, 12).Value2 = 14000
$ws.Cells.Item(102, 13).Value2 = 13500
$ws.Cells.Item(102, 16).Value2 = 675
$ws.Cells.Item(103, 4).Value2 = 44557
$ws.Cells.Item(103, 11).Value2 = 11000
$ws.Cells.Item(103, 12).Value2 = 12000
$ws.Cells.Item(103, 13).Value2 = 11500
$ws.Cells.Item(103, 16).Value2 = 575
$ws.Cells.Item(104, 4).Value2 = 44711
$ws.Cells.Item(104, 10).Value2 = 120
$ws.Cells.Item(104, 11).Value2 = 54000
$ws.Cells.Item(104, 12).Value2 = 55000
$ws.Cells.Item(104, 13).Value2 = 54500
$ws.Cells.Item(104, 16).Value2 = 2725
$ws.Cells.Item(105, 4).Value2 = 44711
$ws.Cells.Item(105, 10).Value2 = 120
$ws.Cells.Item(105, 11).Value2 = 48000
$ws.Cells.Item(105, 12).Value2 = 50000
$ws.Cells.Item(105, 13).Value2 = 49000
$ws.Cells.Item(105, 16).Value2 = 2450
$ws.Cells.Item(106, 4).Value2 = 44260
$ws.Cells.Item(106, 10).Value2 = 160
$ws.Cells.Item(107, 4).Value2 = 44778
$ws.Cells.Item(107, 10).Value2 = 130
$ws.Cells.Item(107, 11).Value2 = 19000
$ws.Cells.Item(107, 12).Value2 = 20000
$ws.Cells.Item(107, 13).Value2 = 19500
$ws.Cells.Item(107, 16).Value2 = 975
$ws.Cells.Item(108, 4).Value2 = 44449
$ws.Cells.Item(108, 11).Value2 = 24000
$ws.Cells.Item(108, 12).Value2 = 25000
$ws.Cells.Item(108, 13).Value2 = 24500
$ws.Cells.Item(108, 16).Value2 = 1225
$ws.Cells.Item(109, 4).Value2 = 44449
$ws.Cells.Item(109, 9).Value2 = "Segunda"
$ws.Cells.Item(109, 10).Value2 = 160
$ws.Cells.Item(109, 11).Value2 = 22000
$ws.Cells.Item(109, 12).Value2 = 23000
$ws.Cells.Item(109, 13).Value2 = 22500
$ws.Cells.Item(109, 16).Value2 = 1125
$ws.Cells.Item(110, 4).Value2 = 44434
$ws.Cells.Item(110, 10).Value2 = 120
$ws.Cells.Item(110, 11).Value2 = 20000
$ws.Cells.Item(110, 12).Value2 = 21000
$ws.Cells.Item(110, 13).Value2 = 20500
$ws.Cells.Item(110, 16).Value2 = 1025
$ws.Cells.Item(111, 4).Value2 = 44442
$ws.Cells.Item(111, 9).Value2 = "Primera"
$ws.Cells.Item(111, 10).Value2 = 120
$ws.Cells.Item(111, 11).Value2 = 15000
$ws.Cells.Item(111, 12).Value2 = 16000
$ws.Cells.Item(111, 13).Value2 = 15500
$ws.Cells.Item(111, 16).Value2 = 775
$ws.Cells.Item(112, 4).Value2 = 44533
$ws.Cells.Item(112, 10).Value2 = 160
$ws.Cells.Item(112, 11).Value2 = 13000
$ws.Cells.Item(112, 12).Value2 = 14000
$ws.Cells.Item(112, 13).Value2 = 13500
$ws.Cells.Item(112, 16).Value2 = 675
$ws.Cells.Item(113, 4).Value2 = 44533
$ws.Cells.Item(113, 10).Value2 = 140
$ws.Cells.Item(113, 11).Value2 = 11000
$ws.Cells.Item(113, 12).Value2 = 12000
$ws.Cells.Item(113, 13).Value2 = 11500
$ws.Cells.Item(113, 16).Value2 = 575
$ws.Cells.Item(114, 4).Value2 = 44421
$ws.Cells.Item(114, 11).Value2 = 23000
$ws.Cells.Item(114, 12).Value2 = 24000
$ws.Cells.Item(114, 13).Value2 = 23500
$ws.Cells.Item(114, 16).Value2 = 1175
$ws.Cells.Item(115, 4).Value2 = 44421
$ws.Cells.Item(115, 10).Value2 = 120
$ws.Cells.Item(115, 11).Value2 = 21000
$ws.Cells.Item(115, 13).Value2 = 21500
$ws.Cells.Item(115, 16).Value2 = 1075
$ws.Cells.Item(116, 4).Value2 = 44589
$ws.Cells.Item(116, 10).Value2 = 120
$ws.Cells.Item(116, 11).Value2 = 30000
$ws.Cells.Item(116, 12).Value2 = 32000
$ws.Cells.Item(116, 13).Value2 = 31000
$ws.Cells.Item(116, 16).Value2 = 1550
$ws.Cells.Item(117, 4).Value2 = 44589
$ws.Cells.Item(117, 11).Value2 = 20000
$ws.Cells.Item(117, 12).Value2 = 22000
$ws.Cells.Item(117, 13).Value2 = 21000
$ws.Cells.Item(117, 16).Value2 = 1050
$ws.Cells.Item(118, 4).Value2 = 44309
$ws.Cells.Item(118, 10).Value2 = 140
$ws.Cells.Item(118, 11).Value2 = 19000
$ws.Cells.Item(118, 12).Value2 = 20000
$ws.Cells.Item(118, 13).Value2 = 19500
$ws.Cells.Item(118, 16).Value2 = 975
$ws.Cells.Item(119, 4).Value2 = 44309
$ws.Cells.Item(119, 9).Value2 = "Segunda"
$ws.Cells.Item(119, 10).Value2 = 160
$ws.Cells.Item(119, 11).Value2 = 17000
$ws.Cells.Item(119, 12).Value2 = 18000
$ws.Cells.Item(119, 13).Value2 = 17500
$ws.Cells.Item(119, 16).Value2 = 875
$ws.Cells.Item(120, 4).Value2 = 44771
$ws.Cells.Item(120, 9).Value2 = "Primera"
$ws.Cells.Item(120, 10).Value2 = 150
$ws.Cells.Item(120, 11).Value2 = 20000
$ws.Cells.Item(120, 12).Value2 = 22000
$ws.Cells.Item(121, 4).Value2 = 44379
$ws.Cells.Item(121, 10).Value2 = 180
$ws.Cells.Item(121, 11).Value2 = 22000
$ws.Cells.Item(121, 12).Value2 = 25000
$ws.Cells.Item(121, 13).Value2 = 23333
$ws.Cells.Item(121, 16).Value2 = 1167
$ws.Cells.Item(122, 4).Value2 = 44379
$ws.Cells.Item(122, 10).Value2 = 80
$ws.Cells.Item(122, 11).Value2 = 21000
$ws.Cells.Item(122, 12).Value2 = 21000
$ws.Cells.Item(122, 13).Value2 = 21000
$ws.Cells.Item(122, 16).Value2 = 1050
$ws.Cells.Item(123, 4).Value2 = 44435
$ws.Cells.Item(123, 10).Value2 = 240
$ws.Cells.Item(123, 11).Value2 = 16000
$ws.Cells.Item(123, 12).Value2 = 21000
$ws.Cells.Item(123, 13).Value2 = 18500
$ws.Cells.Item(123, 16).Value2 = 925
$ws.Cells.Item(124, 4).Value2 = 44435
$ws.Cells.Item(124, 11).Value2 = 14000
$ws.Cells.Item(124, 12).Value2 = 15000
$ws.Cells.Item(124, 13).Value2 = 14500
$ws.Cells.Item(124, 16).Value2 = 725
$ws.Cells.Item(125, 4).Value2 = 44319
$ws.Cells.Item(125, 11).Value2 = 28000
$ws.Cells.Item(125, 12).Value2 = 30000
$ws.Cells.Item(125, 13).Value2 = 29000
$ws.Cells.Item(125, 16).Value2 = 1450
$ws.Cells.Item(126, 4).Value2 = 44319
$ws.Cells.Item(126, 9).Value2 = "Segunda"
$ws.Cells.Item(126, 10).Value2 = 120
$ws.Cells.Item(126, 11).Value2 = 25000
$ws.Cells.Item(126, 12).Value2 = 26000
$ws.Cells.Item(126, 13).Value2 = 25500
$ws.Cells.Item(126, 16).Value2 = 1275
$ws.Cells.Item(127, 4).Value2 = 44344
$ws.Cells.Item(127, 10).Value2 = 120
$ws.Cells.Item(127, 11).Value2 = 45000
$ws.Cells.Item(127, 12).Value2 = 46000
$ws.Cells.Item(127, 13).Value2 = 45500
$ws.Cells.Item(127, 16).Value2 = 2275
$ws.Cells.Item(128, 4).Value2 = 44232
$ws.Cells.Item(128, 9).Value2 = "Primera"
$ws.Cells.Item(128, 10).Value2 = 140
$ws.Cells.Item(128, 11).Value2 = 29000
$ws.Cells.Item(128, 12).Value2 = 30000
$ws.Cells.Item(128, 13).Value2 = 29500
$ws.Cells.Item(128, 16).Value2 = 1475
$ws.Cells.Item(129, 4).Value2 = 44484
$ws.Cells.Item(129, 10).Value2 = 140
$ws.Cells.Item(129, 11).Value2 = 26000
$ws.Cells.Item(129, 12).Value2 = 27000
$ws.Cells.Item(129, 13).Value2 = 26500
$ws.Cells.Item(129, 16).Value2 = 1325
$ws.Cells.Item(130, 4).Value2 = 44484
$ws.Cells.Item(130, 11).Value2 = 24000
$ws.Cells.Item(130, 12).Value2 = 25000
$ws.Cells.Item(130, 13).Value2 = 24500
$ws.Cells.Item(130, 16).Value2 = 1225
$ws.Cells.Item(131, 4).Value2 = 44452
$ws.Cells.Item(131, 10).Value2 = 120
$ws.Cells.Item(131, 11).Value2 = 25000
$ws.Cells.Item(131, 12).Value2 = 26000
$ws.Cells.Item(131, 13).Value2 = 25500
$ws.Cells.Item(131, 16).Value2 = 1275
$ws.Cells.Item(132, 4).Value2 = 44452
$ws.Cells.Item(132, 10).Value2 = 120
$ws.Cells.Item(132, 11).Value2 = 22000
$ws.Cells.Item(132, 12).Value2 = 23000
$ws.Cells.Item(132, 13).Value2 = 22500
$ws.Cells.Item(132, 16).Value2 = 1125

# --- Append new rows 133 and 134 (full records), cloning row 131/132 formatting ---
$ws.Cells.Item(133, 1).Value2 = 1
$ws.Cells.Item(133, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(133, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(133, 4).Value2 = 44701
$ws.Cells.Item(133, 5).Value2 = 15
$ws.Cells.Item(133, 6).Value2 = 100112042
$ws.Cells.Item(133, 7).Value2 = "Locoto"
$ws.Cells.Item(133, 8).Value2 = "Sin especificar"
$ws.Cells.Item(133, 9).Value2 = "Primera"
$ws.Cells.Item(133, 10).Value2 = 100
$ws.Cells.Item(133, 11).Value2 = 54000
$ws.Cells.Item(133, 12).Value2 = 55000
$ws.Cells.Item(133, 13).Value2 = 54500
$ws.Cells.Item(133, 14).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(133, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(133, 16).Value2 = 2725
$ws.Cells.Item(133, 17).Value2 = 20
$ws.Cells.Item(133, 18).Value2 = "Hortaliza"
$ws.Cells.Item(133, 4).NumberFormat = $ws.Cells.Item(131, 4).NumberFormat

$ws.Cells.Item(134, 1).Value2 = 1
$ws.Cells.Item(134, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(134, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(134, 4).Value2 = 44701
$ws.Cells.Item(134, 5).Value2 = 15
$ws.Cells.Item(134, 6).Value2 = 100112042
$ws.Cells.Item(134, 7).Value2 = "Locoto"
$ws.Cells.Item(134, 8).Value2 = "Sin especificar"
$ws.Cells.Item(134, 9).Value2 = "Segunda"
$ws.Cells.Item(134, 10).Value2 = 160
$ws.Cells.Item(134, 11).Value2 = 49000
$ws.Cells.Item(134, 12).Value2 = 50000
$ws.Cells.Item(134, 13).Value2 = 49500
$ws.Cells.Item(134, 14).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(134, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(134, 16).Value2 = 2475
$ws.Cells.Item(134, 17).Value2 = 20
$ws.Cells.Item(134, 18).Value2 = "Hortaliza"
$ws.Cells.Item(134, 4).NumberFormat = $ws.Cells.Item(132, 4).NumberFormat

